$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")

$cell = $ws.Cells.Item(2, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-15'
$cell.Style = "Normal"
$ws.Cells.Item(2, 3).Value = '上饶·次元重现夏日嘉年华（取消）'
$ws.Cells.Item(2, 4).Value = '普济巷地委大院北侧约90米 四季体育运动馆'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.15 09:30-08.15 17:30'
$cell.Style = "Normal"
$ws.Cells.Item(2, 6).Value = 128
$ws.Cells.Item(2, 7).Value = '不可售'
$ws.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87679'
$ws.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'

$cell = $ws.Cells.Item(3, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-15'
$cell.Style = "Normal"
$ws.Cells.Item(3, 3).Value = '乐平·CY境界次元第三届动漫游戏庆典（取消）'
$ws.Cells.Item(3, 4).Value = '乐平大道5号 东方国际酒店'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.15 09:00-08.15 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(3, 6).Value = 77
$ws.Cells.Item(3, 7).Value = '不可售'
$ws.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90005'
$ws.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/VEUop3K21722251709547.png'

$cell = $ws.Cells.Item(4, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-17'
$cell.Style = "Normal"
$ws.Cells.Item(4, 3).Value = '南昌·CM03·配音演员孙路路专场见面会'
$ws.Cells.Item(4, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.17 09:00-08.17 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(4, 6).Value = 116
$ws.Cells.Item(4, 7).Value = 258
$ws.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89623'
$ws.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/TIDNWGni1721565446649.jpeg'

$cell = $ws.Cells.Item(5, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-17'
$cell.Style = "Normal"
$ws.Cells.Item(5, 3).Value = '南昌·CM03动漫游戏博览会'
$ws.Cells.Item(5, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.17 09:00-08.18 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(5, 6).Value = 5372
$ws.Cells.Item(5, 7).Value = 65
$ws.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89295'
$ws.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/YhHLfv5y1722849043508.jpeg'

$cell = $ws.Cells.Item(6, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-18'
$cell.Style = "Normal"
$ws.Cells.Item(6, 3).Value = '九江·如梦令国潮动漫节'
$ws.Cells.Item(6, 4).Value = '十里大道202号 山水国际大酒店(九江火车站快乐城店)'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.18 11:00-08.18 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(6, 6).Value = 71
$ws.Cells.Item(6, 7).Value = 40
$ws.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90126'
$ws.Cells.Item(6, 9).Value = '//i1.hdslb.com/bfs/openplatform/202407/bs3xfiQc1721988224155.jpeg'

$cell = $ws.Cells.Item(7, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-24'
$cell.Style = "Normal"
$ws.Cells.Item(7, 3).Value = '南昌·第四届龙年动漫展——暑假最后的狂欢'
$ws.Cells.Item(7, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.24 10:00-08.25 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(7, 6).Value = 886
$ws.Cells.Item(7, 7).Value = 55
$ws.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87135'
$ws.Cells.Item(7, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'

$cell = $ws.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-24'
$cell.Style = "Normal"
$ws.Cells.Item(8, 3).Value = '赣州·第五人格only'
$ws.Cells.Item(8, 4).Value = '兴国路恒大帝景西门 江西长庚控股有限公司'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.24 10:00-08.24 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(8, 6).Value = 132
$ws.Cells.Item(8, 7).Value = 55
$ws.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89742'
$ws.Cells.Item(8, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/Jxx8Wz6I1721644479535.jpeg'

$cell = $ws.Cells.Item(9, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-15'
$cell.Style = "Normal"
$ws.Cells.Item(9, 3).Value = '南昌·Sunflower Garden动漫游戏展'
$ws.Cells.Item(9, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.15 09:00-09.16 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(9, 6).Value = 2381
$ws.Cells.Item(9, 7).Value = 65
$ws.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89659'
$ws.Cells.Item(9, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/CQCXbg291721632431682.jpeg'

$cell = $ws.Cells.Item(10, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-15'
$cell.Style = "Normal"
$ws.Cells.Item(10, 3).Value = '南昌·第一届哥布林动漫游戏展——开学季&贺中秋'
$ws.Cells.Item(10, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.15 10:00-09.16 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(10, 6).Value = 78
$ws.Cells.Item(10, 7).Value = 55
$ws.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89240'
$ws.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/pixnzm5p1720496832036.jpeg'

$cell = $ws.Cells.Item(11, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-17'
$cell.Style = "Normal"
$ws.Cells.Item(11, 3).Value = '南昌·Aud中秋动漫嘉年华'
$ws.Cells.Item(11, 4).Value = '青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.17 10:00-09.17 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(11, 6).Value = 53
$ws.Cells.Item(11, 7).Value = 29.9
$ws.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90329'
$ws.Cells.Item(11, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/pbU7Eftp1722660514298.jpeg'

$cell = $ws.Cells.Item(12, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-10-02'
$cell.Style = "Normal"
$ws.Cells.Item(12, 3).Value = '南昌·萌卡动漫展'
$ws.Cells.Item(12, 4).Value = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.10.02 09:00-10.03 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(12, 6).Value = 2232
$ws.Cells.Item(12, 7).Value = 65
$ws.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89738'
$ws.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/uqTvacSV1721621530709.jpeg'

$ws.Rows(13).Delete()

$ws = $wb.Worksheets.Item("全部类型")

$cell = $ws.Cells.Item(2, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-15'
$cell.Style = "Normal"
$ws.Cells.Item(2, 3).Value = '上饶·次元重现夏日嘉年华（取消）'
$ws.Cells.Item(2, 4).Value = '普济巷地委大院北侧约90米 四季体育运动馆'
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.15 09:30-08.15 17:30'
$cell.Style = "Normal"
$ws.Cells.Item(2, 6).Value = 128
$ws.Cells.Item(2, 7).Value = '不可售'
$ws.Cells.Item(2, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87679'
$ws.Cells.Item(2, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'

$cell = $ws.Cells.Item(3, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-15'
$cell.Style = "Normal"
$ws.Cells.Item(3, 3).Value = '乐平·CY境界次元第三届动漫游戏庆典（取消）'
$ws.Cells.Item(3, 4).Value = '乐平大道5号 东方国际酒店'
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.15 09:00-08.15 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(3, 6).Value = 77
$ws.Cells.Item(3, 7).Value = '不可售'
$ws.Cells.Item(3, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90005'
$ws.Cells.Item(3, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/VEUop3K21722251709547.png'

$cell = $ws.Cells.Item(4, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-17'
$cell.Style = "Normal"
$ws.Cells.Item(4, 3).Value = '南昌·CM03·配音演员孙路路专场见面会'
$ws.Cells.Item(4, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.17 09:00-08.17 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(4, 6).Value = 116
$ws.Cells.Item(4, 7).Value = 258
$ws.Cells.Item(4, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89623'
$ws.Cells.Item(4, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/TIDNWGni1721565446649.jpeg'

$cell = $ws.Cells.Item(5, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-17'
$cell.Style = "Normal"
$ws.Cells.Item(5, 3).Value = '南昌·CM03动漫游戏博览会'
$ws.Cells.Item(5, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.17 09:00-08.18 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(5, 6).Value = 5372
$ws.Cells.Item(5, 7).Value = 65
$ws.Cells.Item(5, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89295'
$ws.Cells.Item(5, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/YhHLfv5y1722849043508.jpeg'

$cell = $ws.Cells.Item(6, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-17'
$cell.Style = "Normal"
$ws.Cells.Item(6, 3).Value = '南昌·CrossingX意次元｜乐队番ONLY'
$ws.Cells.Item(6, 4).Value = '佘山路66号樟树林文化生活公园2座 黑铁Livehouse'
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.17 14:00-08.17 15:30'
$cell.Style = "Normal"
$ws.Cells.Item(6, 6).Value = 94
$ws.Cells.Item(6, 7).Value = 38
$ws.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89670'
$ws.Cells.Item(6, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/t07f8mmz1721894837940.jpeg'

$cell = $ws.Cells.Item(7, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-18'
$cell.Style = "Normal"
$ws.Cells.Item(7, 3).Value = '九江·如梦令国潮动漫节'
$ws.Cells.Item(7, 4).Value = '十里大道202号 山水国际大酒店(九江火车站快乐城店)'
$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.18 11:00-08.18 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(7, 6).Value = 71
$ws.Cells.Item(7, 7).Value = 40
$ws.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90126'
$ws.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202407/bs3xfiQc1721988224155.jpeg'

$cell = $ws.Cells.Item(8, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-24'
$cell.Style = "Normal"
$ws.Cells.Item(8, 3).Value = '南昌·【8月24日】滑稽互动狂欢大作战《欢乐小丑嘉年华》'
$ws.Cells.Item(8, 4).Value = '象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院'
$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.24 14:30-08.24 20:00'
$cell.Style = "Normal"
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 60
$ws.Cells.Item(8, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90177'
$ws.Cells.Item(8, 9).Value = '//i1.hdslb.com/bfs/openplatform/202408/ed1EQGH71722479594577.jpeg'

$cell = $ws.Cells.Item(9, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-24'
$cell.Style = "Normal"
$ws.Cells.Item(9, 3).Value = '南昌·第四届龙年动漫展——暑假最后的狂欢'
$ws.Cells.Item(9, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.24 10:00-08.25 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(9, 6).Value = 886
$ws.Cells.Item(9, 7).Value = 55
$ws.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87135'
$ws.Cells.Item(9, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'

$cell = $ws.Cells.Item(10, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-08-24'
$cell.Style = "Normal"
$ws.Cells.Item(10, 3).Value = '赣州·第五人格only'
$ws.Cells.Item(10, 4).Value = '兴国路恒大帝景西门 江西长庚控股有限公司'
$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.08.24 10:00-08.24 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(10, 6).Value = 132
$ws.Cells.Item(10, 7).Value = 55
$ws.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89742'
$ws.Cells.Item(10, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/Jxx8Wz6I1721644479535.jpeg'

$cell = $ws.Cells.Item(11, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-15'
$cell.Style = "Normal"
$ws.Cells.Item(11, 3).Value = '南昌·Sunflower Garden动漫游戏展'
$ws.Cells.Item(11, 4).Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.15 09:00-09.16 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(11, 6).Value = 2381
$ws.Cells.Item(11, 7).Value = 65
$ws.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89659'
$ws.Cells.Item(11, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/CQCXbg291721632431682.jpeg'

$cell = $ws.Cells.Item(12, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-15'
$cell.Style = "Normal"
$ws.Cells.Item(12, 3).Value = '南昌·第一届哥布林动漫游戏展——开学季&贺中秋'
$ws.Cells.Item(12, 4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.15 10:00-09.16 18:00'
$cell.Style = "Normal"
$ws.Cells.Item(12, 6).Value = 78
$ws.Cells.Item(12, 7).Value = 55
$ws.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89240'
$ws.Cells.Item(12, 9).Value = '//i2.hdslb.com/bfs/openplatform/202407/pixnzm5p1720496832036.jpeg'

$cell = $ws.Cells.Item(13, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-16'
$cell.Style = "Normal"
$ws.Cells.Item(13, 3).Value = '南昌·《梁祝》65周年大型交响音乐会-风兔子交响乐团'
$ws.Cells.Item(13, 4).Value = '象湖新城东祥路昌南文化中心五号馆 昌南文化中心大剧院'
$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.16 19:30-09.16 20:50'
$cell.Style = "Normal"
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 140
$ws.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90515'
$ws.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202408/muKn0Ygv1723107475651.jpeg'

$cell = $ws.Cells.Item(14, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-09-17'
$cell.Style = "Normal"
$ws.Cells.Item(14, 3).Value = '南昌·Aud中秋动漫嘉年华'
$ws.Cells.Item(14, 4).Value = '青山湖南大道260号泰丰轮胎厂进大门走到底左拐 赣A篮球梦时代GANA PARK'
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.09.17 10:00-09.17 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(14, 6).Value = 53
$ws.Cells.Item(14, 7).Value = 29.9
$ws.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=90329'
$ws.Cells.Item(14, 9).Value = '//i0.hdslb.com/bfs/openplatform/202408/pbU7Eftp1722660514298.jpeg'

$cell = $ws.Cells.Item(15, 2)
$cell.NumberFormat = "@"
$cell.Value = '2024-10-02'
$cell.Style = "Normal"
$ws.Cells.Item(15, 3).Value = '南昌·萌卡动漫展'
$ws.Cells.Item(15, 4).Value = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = '2024.10.02 09:00-10.03 17:00'
$cell.Style = "Normal"
$ws.Cells.Item(15, 6).Value = 2232
$ws.Cells.Item(15, 7).Value = 65
$ws.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=89738'
$ws.Cells.Item(15, 9).Value = '//i0.hdslb.com/bfs/openplatform/202407/uqTvacSV1721621530709.jpeg'

$ws.Rows(16).Delete()
